$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove all existing merges (the target layout has no merged cells) ---
$ws.Cells.UnMerge()

# --- Clear out the whole former table area (A1:D14) so we can rebuild it cleanly ---
$ws.Range("A1:D14").Clear()

# --- Row 1 ---
$ws.Range("C1").Value = "Barracks"
$ws.Range("C1").Font.Size = 8
$ws.Range("D1").Value = ""
$ws.Range("D1").Font.Size = 8

# --- Row 2 ---
$ws.Range("A2").Value = "Vault_Door"
$ws.Range("A2").Font.Size = 8
$ws.Range("B2").Value = "Elevator(0)"
$ws.Range("B2").Font.Size = 8
$ws.Range("C2").Value = "Corridor1"
$ws.Range("C2").Font.Size = 8
$ws.Range("D2").Value = ""
$ws.Range("D2").Font.Size = 8

# --- Row 3 ---
$ws.Range("B3").Value = "."
$ws.Range("B3").Font.Size = 8
$ws.Range("C3").Value = "Dining_Room"
$ws.Range("C3").Font.Size = 8

# --- Row 4 ---
$ws.Range("B4").Value = "."
$ws.Range("B4").Font.Size = 8

# --- Row 5 ---
$ws.Range("B5").Value = "."
$ws.Range("B5").Font.Size = 8
$ws.Range("C5").Value = "Water_Treatment_Plant"
$ws.Range("C5").Font.Size = 8
$ws.Range("D5").Value = ""
$ws.Range("D5").Font.Size = 8

# --- Row 6 ---
$ws.Range("A6").Value = "Storage_Room"
$ws.Range("A6").Font.Size = 8
$ws.Range("B6").Value = "Elevator(-1)"
$ws.Range("B6").Font.Size = 8
$ws.Range("C6").Value = "Corridor2"
$ws.Range("C6").Font.Size = 8
$ws.Range("D6").Value = ""
$ws.Range("D6").Font.Size = 8

# --- Row 7 ---
$ws.Range("B7").Value = "."
$ws.Range("B7").Font.Size = 8
$ws.Range("C7").Value = "Power_Plant "
$ws.Range("C7").Font.Size = 8

# --- Row 8 ---
$ws.Range("B8").Value = "."
$ws.Range("B8").Font.Size = 8

# --- Row 9 ---
$ws.Range("A9").Value = "Armory"
$ws.Range("A9").Font.Size = 8
$ws.Range("B9").Value = "."
$ws.Range("B9").Font.Size = 8
$ws.Range("C9").Value = "Hospital "
$ws.Range("C9").Font.Size = 8

# --- Row 10 ---
$ws.Range("A10").Value = "Corridor4"
$ws.Range("A10").Font.Size = 8
$ws.Range("B10").Value = "Elevator(-2)"
$ws.Range("B10").Font.Size = 8
$ws.Range("C10").Value = "Corridor3"
$ws.Range("C10").Font.Size = 8

# --- Row 11 ---
$ws.Range("A11").Value = "weight_room"
$ws.Range("A11").Font.Size = 8
$ws.Range("C11").Value = "Science_Lab"
$ws.Range("C11").Font.Size = 8
$ws.Range("D11").Value = ""
$ws.Range("D11").Font.Size = 8

# --- Update selection to match the saved view state ---
[void]$ws.Range("B2").Select()
